$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" column (D) and "Volume(1h)" column (E) values.
# Cells whose new value parses as a plain number need their number format
# forced to text first, otherwise Excel will auto-convert the text into a
# numeric value (losing formatting like trailing zeros or leading zeros).

$ws.Range("D2").Value = '44.741.46'
$ws.Range("E2").Value = '  +4.30%  '

$ws.Range("D3").Value = '2.420.24'
$ws.Range("E3").Value = '  +2.57%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.53'
$ws.Range("E5").Value = '  +4.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.34'
$ws.Range("E6").Value = '  +6.55%  '

$ws.Range("E7").Value = '  +2.20%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("E9").Value = '  +11.17%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.34'
$ws.Range("E10").Value = '  +3.19%  '

$ws.Range("E11").Value = '  +1.81%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.71'
$ws.Range("E12").Value = '  +1.81%  '

$ws.Range("E13").Value = '  -2.21%  '

$ws.Range("E14").Value = '  +2.94%  '

$ws.Range("D15").Value = '2.799.32'
$ws.Range("E15").Value = '  +2.70%  '

$ws.Range("D16").Value = '2.416.49'
$ws.Range("E16").Value = '  +1.64%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.831'
$ws.Range("E17").Value = '  +4.72%  '

$ws.Range("D18").Value = '44.518.22'
$ws.Range("E18").Value = '  +3.72%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.34'
$ws.Range("E19").Value = '  +3.28%  '

$ws.Range("E20").Value = '  +1.77%  '

$ws.Range("E21").Value = '  +3.65%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.73'
$ws.Range("E22").Value = '  +1.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '242.23'
$ws.Range("E23").Value = '  +3.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.27'
$ws.Range("E24").Value = '  +4.11%  '

$ws.Range("E25").Value = '  +2.32%  '

$ws.Range("E26").Value = '  -0.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.17'
$ws.Range("E27").Value = '  +3.49%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.29'
$ws.Range("E28").Value = '  -3.03%  '

$ws.Range("E29").Value = '  +2.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.60'
$ws.Range("E30").Value = '  +3.75%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '48.50'
$ws.Range("E31").Value = '  +1.34%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.127'
$ws.Range("E32").Value = '  +18.28%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.47'
$ws.Range("E33").Value = '  +11.35%  '

$ws.Range("E34").Value = '  +3.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0771'
$ws.Range("E35").Value = '  +6.47%  '

$ws.Range("E36").Value = '  +0.27%  '

$ws.Range("E37").Value = '  +3.74%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.49'
$ws.Range("E38").Value = '  +3.96%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.84'
$ws.Range("E39").Value = '  +0.12%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '123.35'
$ws.Range("E40").Value = '  -4.60%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.110'
$ws.Range("E41").Value = '  +1.98%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.20'
$ws.Range("E42").Value = '  -3.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.02'
$ws.Range("E43").Value = '  +0.26%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0290'
$ws.Range("E44").Value = '  +4.36%  '

$ws.Range("D45").Value = '1.939.31'
$ws.Range("E45").Value = '  +0.66%  '

$ws.Range("E46").Value = '  -1.50%  '

$ws.Range("E47").Value = '  +8.12%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.26'
$ws.Range("E48").Value = '  -0.19%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.73'
$ws.Range("E49").Value = '  +15.33%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.72'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.93'
$ws.Range("E51").Value = '  +5.56%  '
